# Add a new "Erp Code" column (E) to the PLACE sheet, populate it for the
# existing 95 data rows, append 3 new data rows (96-98), resize columns,
# and refresh the window view (zoom + selection) - matches the commit's
# "Some bug fixed and upgraded" change that tags each placement with its
# ERP part code and three extra test points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -----------------------------------------------------
$ws.Range("E1").Value = "Erp Code"

# --- Erp Code values for existing rows 2-48 (the same 47-row pattern repeats
#     for rows 49-95, mirroring the existing A:D duplication in the sheet) --
$ErpCodes = @(
    130100009,130100001,130100010,130100014,130100009,130100009,130100011,
    130100010,130100010,130100009,130100009,130100009,130100009,130100009,
    130100009,130100009,130100009,130100014,130100016,130100016,130100010,
    130100010,130100010,130100010,130100010,130100010,130100010,130100010,
    130100010,130100011,130100010,130100011,130100011,130100010,130100011,
    130100011,130100010,130100011,130100011,130100015,130100012,130100012,
    130100012,130100012,130100012,130100014,130100014
)

for ($i = 0; $i -lt $ErpCodes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $ErpCodes[$i]
}
for ($i = 0; $i -lt $ErpCodes.Length; $i++) {
    $row = $i + 49
    $ws.Cells.Item($row, 5).Value = $ErpCodes[$i]
}

# --- Three brand-new test-point rows (96-98) --------------------------------
$ws.Range("A96").Value = "99.9999mm"
$ws.Range("B96").Value = "55.5555mm"
$ws.Range("C96").Value = 90
$ws.Range("D96").Value = "T1"
$ws.Range("E96").Value = 130100099

$ws.Range("A97").Value = "999.9999mm"
$ws.Range("B97").Value = "555.5555mm"
$ws.Range("C97").Value = 90
$ws.Range("D97").Value = "T2"
$ws.Range("E97").Value = 130100098

$ws.Range("A98").Value = "99.9999mm"
$ws.Range("B98").Value = "55.5555mm"
$ws.Range("C98").Value = 90
$ws.Range("D98").Value = "T3"
$ws.Range("E98").Value = 130100097

# --- Give the whole new column E the same centered format as column D ------
$ws.Range("D1:D98").Copy()
$ws.Range("E1:E98").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = 0

# --- Column E width (matches the new <col> entry in the diff) --------------
# NOTE: ColumnWidth is quantized internally to 1/6-character steps, so the
# nearest representable value to the target OOXML width (13.140625) is used.
$ws.Columns.Item(5).ColumnWidth = 12.333333333333334

# --- View refresh: 130% zoom + active selection on K32 ----------------------
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("K32").Select()

Write-Output "done"
